$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table rows 16..69 (column E = "Periodo Mora" code, column F = "Valor Mora")
# were reordered: the block is reversed top-to-bottom (row 16 <-> row 69,
# row 17 <-> row 68, ... row 42 <-> row 43), swapping both the period label
# (E) and the mora value (F) together so each period keeps its own value.
$firstRow = 16
$lastRow = 69
$pairCount = [Math]::Floor((($lastRow - $firstRow + 1)) / 2)

for ($i = 0; $i -lt $pairCount; $i++) {
    $rowA = $firstRow + $i
    $rowB = $lastRow - $i

    $eA = $ws.Cells.Item($rowA, 5).Value2
    $fA = $ws.Cells.Item($rowA, 6).Value2
    $eB = $ws.Cells.Item($rowB, 5).Value2
    $fB = $ws.Cells.Item($rowB, 6).Value2

    $ws.Cells.Item($rowA, 5).Value = $eB
    $ws.Cells.Item($rowA, 6).Value = $fB
    $ws.Cells.Item($rowB, 5).Value = $eA
    $ws.Cells.Item($rowB, 6).Value = $fA
}
